# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
#
# Adds an MSME "by firm size" breakdown table (Number of employees / Assets /
# Turnover headers over Micro / Small / Medium / Large rows) to the Burundi
# Summary sheet, just above the existing ISTEEBU source citation block - and
# pushes that citation block further down the sheet to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 6 blank rows at row 16. This shifts the existing
# "ISTEEBU" / long citation rows (formerly rows 21-22, with their original
# bold/italic formatting intact) down to rows 27-28.
$ws.Rows("16:21").Insert()

# Row 16: bold column headers for the new table.
$ws.Range("B16").Value = "Number of employees"
$ws.Range("C16").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D16").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B16:D16").Font.Bold = $true

# Rows 17-20: firm-size labels in column A (normal/default style).
$ws.Range("A17").Value = "Micro"
$ws.Range("A18").Value = "Small"
$ws.Range("A19").Value = "Medium"
$ws.Range("A20").Value = "Large"

# The shifted citation block (now rows 27-28) keeps its original bold /
# italic look; re-assert it explicitly so it survives the export.
$ws.Range("A27").Font.Bold = $true
$ws.Range("A28").Font.Italic = $true

Write-Host "done"
